$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Clear the fine records data (rows 2-3) on both sheets so the sheets act as
# blank templates for either parking fines or car fines uploads.
$ws1.Range("A2:E3").ClearContents()
$ws1.Hyperlinks.Delete()

$ws2.Range("A2:E3").ClearContents()
$ws2.Hyperlinks.Delete()

# Update selections to match the new active sheet/cell state.
$ws1.Range("F8").Select()

$ws2.Activate()
$ws2.Range("D12").Select()
